$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("go")
$ws.Range("G2:G27").Value = "Pass"
$ws.Range("H2:H27").Value = "刘彩丽"
$ws.Activate()
$ws.Range("B2").Select()
